# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates across sheets per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1859.6
$ws.Range("I43").Value = 1799.8572
$ws.Range("J43").Value = 1911.875
$ws.Range("K43").Value = 1799.8572
$ws.Range("L43").Value = 1911.875
$ws.Range("M43").Value = -1730.8572
$ws.Range("N43").Value = -2049.875
$ws.Range("H58").Value = 4000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 12000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -12300
$ws.Range("H87").Value = 10982.361
$ws.Range("I87").Value = 3410.5
$ws.Range("K87").Value = 3410.5
$ws.Range("M87").Value = -2162.5
$ws.Range("H90").Value = 10982.361
$ws.Range("I90").Value = 3410.5
$ws.Range("K90").Value = 10231.5
$ws.Range("M90").Value = -3991.5
$ws.Range("H98").Value = 6543973.5
$ws.Range("I98").Value = 8527.4375
$ws.Range("J98").Value = 111111110
$ws.Range("K98").Value = 8527.4375
$ws.Range("L98").Value = 111111110
$ws.Range("M98").Value = -7029.4375
$ws.Range("N98").Value = -111114106
$ws.Range("H100").Value = 2203.1667
$ws.Range("I100").Value = 1448.8235
$ws.Range("J100").Value = 3189.6155
$ws.Range("K100").Value = 1448.8235
$ws.Range("L100").Value = 3189.6155
$ws.Range("M100").Value = -907.8235
$ws.Range("N100").Value = -4271.6155
$ws.Range("H112").Value = 6825.9644
$ws.Range("I112").Value = 896.6667
$ws.Range("J112").Value = 7537.48
$ws.Range("K112").Value = 2690.0001
$ws.Range("L112").Value = 22612.44
$ws.Range("M112").Value = -1582.0001
$ws.Range("N112").Value = -24828.44
$ws.Range("H122").Value = 6543973.5
$ws.Range("I122").Value = 8527.4375
$ws.Range("J122").Value = 111111110
$ws.Range("K122").Value = 25582.3125
$ws.Range("L122").Value = 333333330
$ws.Range("M122").Value = -23132.3125
$ws.Range("N122").Value = -333338230
$ws.Range("H139").Value = 70520
$ws.Range("J139").Value = 70520
$ws.Range("L139").Value = 70520
$ws.Range("N139").Value = -80800
$ws.Range("H140").Value = 71402.38
$ws.Range("J140").Value = 91963.336
$ws.Range("L140").Value = 91963.336
$ws.Range("N140").Value = -102323.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2858.5
$ws.Range("I2").Value = 1695.3
$ws.Range("J2").Value = 4312.5
$ws.Range("K2").Value = 1695.3
$ws.Range("L2").Value = 4312.5
$ws.Range("M2").Value = -1582.3
$ws.Range("N2").Value = -4538.5
$ws.Range("H45").Value = 3141.3333
$ws.Range("I45").Value = 3244.8572
$ws.Range("K45").Value = 3244.8572
$ws.Range("M45").Value = -2867.8572
$ws.Range("H61").Value = 2353.639
$ws.Range("I61").Value = 2332.8064
$ws.Range("J61").Value = 2482.8
$ws.Range("K61").Value = 2332.8064
$ws.Range("L61").Value = 2482.8
$ws.Range("M61").Value = -2120.8064
$ws.Range("N61").Value = -2906.8
$ws.Range("H74").Value = 1131.5
$ws.Range("I74").Value = 725.4286
$ws.Range("J74").Value = 1700
$ws.Range("K74").Value = 725.4286
$ws.Range("L74").Value = 1700
$ws.Range("M74").Value = 148.5714
$ws.Range("N74").Value = -3448
$ws.Range("H77").Value = 1131.5
$ws.Range("I77").Value = 725.4286
$ws.Range("J77").Value = 1700
$ws.Range("K77").Value = 3627.143
$ws.Range("L77").Value = 8500
$ws.Range("M77").Value = 740.857
$ws.Range("N77").Value = -17236
$ws.Range("H116").Value = 2858.5
$ws.Range("I116").Value = 1695.3
$ws.Range("J116").Value = 4312.5
$ws.Range("K116").Value = 1695.3
$ws.Range("L116").Value = 4312.5
$ws.Range("M116").Value = 598.7
$ws.Range("N116").Value = -8900.5
$ws.Range("H132").Value = 1728.537
$ws.Range("I132").Value = 1443.1632
$ws.Range("J132").Value = 4525.2
$ws.Range("K132").Value = 4329.4896
$ws.Range("L132").Value = 13575.6
$ws.Range("M132").Value = -1799.4896
$ws.Range("N132").Value = -18635.6
$ws.Range("H136").Value = 2353.639
$ws.Range("I136").Value = 2332.8064
$ws.Range("J136").Value = 2482.8
$ws.Range("K136").Value = 6998.4192
$ws.Range("L136").Value = 7448.400000000001
$ws.Range("M136").Value = -4448.4192
$ws.Range("N136").Value = -12548.4
$ws.Range("H138").Value = 68685.71000000001
$ws.Range("J138").Value = 68685.71000000001
$ws.Range("L138").Value = 68685.71000000001
$ws.Range("N138").Value = -78965.71000000001
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280
$ws.Range("H141").Value = 59928.57
$ws.Range("J141").Value = 62846.152
$ws.Range("L141").Value = 62846.152
$ws.Range("N141").Value = -73206.152

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2858.5
$ws.Range("I3").Value = 1695.3
$ws.Range("J3").Value = 4312.5
$ws.Range("K3").Value = 1695.3
$ws.Range("L3").Value = 4312.5
$ws.Range("M3").Value = -1581.3
$ws.Range("N3").Value = -4540.5
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()
$ws.Range("H140").Value = 59200
$ws.Range("J140").Value = 59200
$ws.Range("L140").Value = 59200
$ws.Range("N140").Value = -69560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 750
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -350
$ws.Range("N22").Value = -1500
$ws.Range("H28").Value = 23500
$ws.Range("J28").Value = 23500
$ws.Range("L28").Value = 23500
$ws.Range("N28").Value = -23990
$ws.Range("H51").Value = 6599.5
$ws.Range("H60").Value = 7899.2
$ws.Range("J60").Value = 8100.75
$ws.Range("L60").Value = 8100.75
$ws.Range("N60").Value = -9122.75
$ws.Range("H61").Value = 6599.5
$ws.Range("H68").Value = 17073
$ws.Range("J68").Value = 17073
$ws.Range("L68").Value = 17073
$ws.Range("N68").Value = -18571
$ws.Range("H71").Value = 17073
$ws.Range("J71").Value = 17073
$ws.Range("L71").Value = 51219
$ws.Range("N71").Value = -58707
$ws.Range("H74").Value = 18735.5
$ws.Range("J74").Value = 18735.5
$ws.Range("L74").Value = 18735.5
$ws.Range("N74").Value = -20483.5
$ws.Range("H77").Value = 18735.5
$ws.Range("J77").Value = 18735.5
$ws.Range("L77").Value = 56206.5
$ws.Range("N77").Value = -64942.5
$ws.Range("H92").Value = 25500
$ws.Range("J92").Value = 25500
$ws.Range("L92").Value = 25500
$ws.Range("N92").Value = -30492
$ws.Range("H138").Value = 47975
$ws.Range("J138").Value = 47975
$ws.Range("L138").Value = 47975
$ws.Range("N138").Value = -58255

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 793.625
$ws.Range("I113").Value = 522.1111
$ws.Range("J113").Value = 821.7126500000001
$ws.Range("K113").Value = 1566.3333
$ws.Range("L113").Value = 2465.13795
$ws.Range("M113").Value = 603.6667000000002
$ws.Range("N113").Value = -6805.13795
$ws.Range("H131").Value = 10000891
$ws.Range("I131").Value = 910.381
$ws.Range("J131").Value = 17242256
$ws.Range("K131").Value = 2731.143
$ws.Range("L131").Value = 51726768
$ws.Range("M131").Value = 2308.857
$ws.Range("N131").Value = -51736848

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2904.3438
$ws.Range("I132").Value = 2368.1738
$ws.Range("J132").Value = 4274.5557
$ws.Range("K132").Value = 7104.5214
$ws.Range("L132").Value = 12823.6671
$ws.Range("M132").Value = -4574.5214
$ws.Range("N132").Value = -17883.6671
$ws.Range("H138").Value = 68750
$ws.Range("J138").Value = 68750
$ws.Range("L138").Value = 68750
$ws.Range("N138").Value = -79030
$ws.Range("H140").Value = 89984.5
$ws.Range("J140").Value = 89984.5
$ws.Range("L140").Value = 89984.5
$ws.Range("N140").Value = -100344.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 658.36365
$ws.Range("I16").Value = 605
$ws.Range("J16").Value = 800.6667
$ws.Range("K16").Value = 605
$ws.Range("L16").Value = 800.6667
$ws.Range("M16").Value = -435
$ws.Range("N16").Value = -1140.6667
$ws.Range("H22").Value = 623.25
$ws.Range("I22").Value = 524.0476
$ws.Range("J22").Value = 762.13336
$ws.Range("K22").Value = 524.0476
$ws.Range("L22").Value = 762.13336
$ws.Range("M22").Value = -229.0476
$ws.Range("N22").Value = -1352.13336
$ws.Range("H27").Value = 623.25
$ws.Range("I27").Value = 524.0476
$ws.Range("J27").Value = 762.13336
$ws.Range("K27").Value = 524.0476
$ws.Range("L27").Value = 762.13336
$ws.Range("M27").Value = -417.0476
$ws.Range("N27").Value = -976.13336
$ws.Range("H136").Value = 2586.8223
$ws.Range("I136").Value = 2751.6843
$ws.Range("J136").Value = 2466.3462
$ws.Range("K136").Value = 8255.052899999999
$ws.Range("L136").Value = 7399.0386
$ws.Range("M136").Value = -5705.052899999999
$ws.Range("N136").Value = -12499.0386
$ws.Range("H138").Value = 53323.832
$ws.Range("J138").Value = 53323.832
$ws.Range("L138").Value = 53323.832
$ws.Range("N138").Value = -63603.832

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1051.75
$ws.Range("I136").Value = 839.3200000000001
$ws.Range("J136").Value = 1534.5454
$ws.Range("K136").Value = 2517.96
$ws.Range("L136").Value = 4603.6362
$ws.Range("M136").Value = 32.03999999999996
$ws.Range("N136").Value = -9703.636200000001
$ws.Range("H139").Value = 61840
$ws.Range("J139").Value = 61840
$ws.Range("L139").Value = 61840
$ws.Range("N139").Value = -72120
